# Dutch translations for HIVE TEAMS.docx

$d = $word.ActiveDocument
$nbsp = [char]160

function Replace-Text {
    param($findText, $replaceText)
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# --- HIVE TEAM: QUALITY ASSURANCE section --------------------------------

# Drop the standalone leading "nbsp" run in front of the heading run.
$target = $nbsp + "HIVE TEAM: QUALITY ASSURANCE"
$rng = $d.Content
$rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$nbspRange = $d.Range($rng.Start, $rng.Start + 1)
$nbspRange.Delete()

$f = "HIVE TEAM: QUALITY ASSURANCE"
$r = "HIVE TEAM: KWALITEITSBORGING"
Replace-Text $f $r

$f = "Ensuring all development tasks meet quality criteria."
$r = "Zorgen dat alle ontwikkeltaken voldoen aan de kwaliteitscriteria."
Replace-Text $f $r

$f = "Release Coordinator"
$r = "Release-coördinator"
Replace-Text $f $r

$f = "auditor"
$r = "controller"
Replace-Text $f $r

$f = "security"
$r = "beveiliging"
Replace-Text $f $r

# --- HIVE TEAM: DEVELOPMENT section ---------------------------------------

$f = $nbsp + "HIVE TEAM: DEVELOPMENT"
$r = "HIVE TEAM: ONTWIKKELING"
Replace-Text $f $r

$f = "Responsible for building SmartCash and supporting applications."
$r = "Verantwoordelijk voor het bouwen van SmartCash en bieden van applicatie ondersteuning."
Replace-Text $f $r

$f = "Hive Coordinator"
$r = "Hive-coördinator"
Replace-Text $f $r

$f = "Creator of the Dash N Drink Soda Machine & SmartCash POS."
$r = "Maker van de Dash N Drink Soda Machine & SmartCash POS."
Replace-Text $f $r

$f = "Developer"
$r = "Ontwikkelaar"
Replace-Text $f $r

$f = "C++ Software Engineer"
$r = "C ++ Software Ontwikkelaar"
Replace-Text $f $r

$f = "Developer"
$r = "Ontwikkelaar"
Replace-Text $f $r

# --- HIVE TEAM: OUTREACH 2 section ----------------------------------------

$f = $nbsp + "HIVE TEAM: OUTREACH 2"
$r = "HIVE TEAM: OUTREACH 2"
Replace-Text $f $r

$f = "This team focuses on community building, growth, general user acquisition in South America"
$r = "Dit team richt zich op gemeenschapsopbouw, groei en het binnen halen van nieuwe gebruikers in Zuid-Amerika"
Replace-Text $f $r

$f = "Hive Coordinator"
$r = "Hive-coördinator"
Replace-Text $f $r

$f = "Outreach Support"
$r = "Outreach ondersteuning"
Replace-Text $f $r

$f = "Outreach Support"
$r = "Outreach ondersteuning"
Replace-Text $f $r

# --- HIVE TEAM: SUPPORT & WEB section --------------------------------------

$f = $nbsp + "HIVE TEAM: SUPPORT" + $nbsp
$r = "HIVE TEAM: ONDERSTEUNING" + $nbsp
Replace-Text $f $r

$f = "This Hive is responsible for on-boarding & generalized SmartCash support."
$r = "Deze Hive is verantwoordelijk voor on-boarding en algemene SmartCash ondersteuning."
Replace-Text $f $r

$f = "Hive Coordinator"
$r = "Hive-coördinator"
Replace-Text $f $r

$f = "Alex is a jack of all trades who loves Technology, Graphics, Web Design & Infrastructure."
$r = "Alex is een manusje-van-alles die houdt van technologie, grafische vormgeving, webdesign en infrastructuur."
Replace-Text $f $r

$f = "Fiscal Officer"
$r = "Fiscaal specialist"
Replace-Text $f $r

$f = "Support"
$r = "Ondersteuning"
Replace-Text $f $r

$f = "Support"
$r = "Ondersteuning"
Replace-Text $f $r

$f = "Assistant Coordinator"
$r = "Assistent-coördinator"
Replace-Text $f $r
